$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Section: createPaypalOrder
$ws.Range("A16").Value = "createPaypalOrder"

$ws.Range("A17").Value = "TestData"
$ws.Range("B17").Value = "intent"
$ws.Range("C17").Value = "currency_code"
$ws.Range("D17").Value = "value"

$ws.Range("A18").Value = "Yes"
$ws.Range("B18").Value = "CAPTURE"
$ws.Range("C18").Value = "INR"
$ws.Range("D18").Value = 500

# Section: getPaypalOrder
$ws.Range("A20").Value = "getPaypalOrder"

$ws.Range("A21").Value = "TestData"
$ws.Range("B21").Value = "intent"
$ws.Range("C21").Value = "currency_code"
$ws.Range("D21").Value = "value"
$ws.Range("E21").Value = "orderId"

$ws.Range("A22").Value = "Yes"
$ws.Range("B22").Value = "CAPTURE"
$ws.Range("C22").Value = "INR"
$ws.Range("D22").Value = 500

# Correct the currency for the createPaypalOrder test to USD
$ws.Range("C18").Value = "USD"

$ws.Range("E22").Value = "7YA93243NU2020819"

$ws.Range("E22").Select()
